$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, "Libre", 2),
    @(2, "Libre", 4),
    @(3, "Libre", 6),
    @(4, "Libre", 3),
    @(5, "Libre", 5),
    @(6, "Libre", 8)
)

$row = 12
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 3).Value = $item[1]
    $ws.Cells.Item($row, 4).Value = $item[2]
    $row = $row + 1
}
